# Update CAN "popis poruka" (message descriptions) on the Translation sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# --- Shorten a couple of existing labels -------------------------------
$ws.Range("F701").Value = "Vol:"
$ws.Range("F704").Value = "Curr:"

# --- Repurpose the old "CELL 5" / "CELL 6" temperature cells -----------
$ws.Range("F751").Value = "TEMP1"
$ws.Range("F752").Value = "TEMP2"

# --- Shorten "SHD cmd" to "SHD" -----------------------------------------
$ws.Range("F765").Value = "SHD"

# --- Append new rows 772-782 with additional CAN message texts ---------
$ws.Range("B772").Value = "SingleUseId884"
$ws.Range("C772").Value = "Typography_02"
$ws.Range("D772").Value = "Left"
$ws.Range("E772").Value = "LTR"
$ws.Range("F772").Value = "TEMP3"

$ws.Range("B773").Value = "SingleUseId885"
$ws.Range("C773").Value = "Typography_02"
$ws.Range("D773").Value = "Left"
$ws.Range("E773").Value = "LTR"
$ws.Range("F773").Value = "TEMP4"

$ws.Range("B774").Value = "SingleUseId886"
$ws.Range("C774").Value = "Typography_02"
$ws.Range("D774").Value = "Left"
$ws.Range("E774").Value = "LTR"
$ws.Range("F774").Value = "TEMP5"

$ws.Range("B775").Value = "SingleUseId887"
$ws.Range("C775").Value = "Typography_05"
$ws.Range("D775").Value = "Left"
$ws.Range("E775").Value = "LTR"
$ws.Range("F775").Value = "Fatal error"

$ws.Range("B776").Value = "SingleUseId888"
$ws.Range("C776").Value = "Typography_02"
$ws.Range("D776").Value = "Left"
$ws.Range("E776").Value = "LTR"
$ws.Range("F776").Value = "TEMP6"

$ws.Range("B777").Value = "SingleUseId889"
$ws.Range("C777").Value = "Typography_01"
$ws.Range("D777").Value = "Left"
$ws.Range("E777").Value = "LTR"
$ws.Range("F777").Value = "SOC:"

$ws.Range("B778").Value = "SingleUseId890"
$ws.Range("C778").Value = "Typography_06"
$ws.Range("D778").Value = "Left"
$ws.Range("E778").Value = "LTR"
$ws.Range("F778").Value = "<value>%"

$ws.Range("B779").Value = "SingleUseId891"
$ws.Range("C779").Value = "Typography_06"
$ws.Range("D779").Value = "Left"
$ws.Range("E779").Value = "LTR"
# "0" would be auto-coerced to a number via .Value (and would also pick up
# a stray number-format style); copy an existing plain-text "0" cell
# instead so the result stays a shared string with the default style.
$ws.Range("F703").Copy() | Out-Null
$ws.Range("F779").PasteSpecial() | Out-Null

$ws.Range("B780").Value = "SingleUseId892"
$ws.Range("C780").Value = "Default"
$ws.Range("D780").Value = "Left"
$ws.Range("E780").Value = "LTR"
$ws.Range("F780").Value = "Actuator fault code:"

$ws.Range("B781").Value = "SingleUseId893"
$ws.Range("C781").Value = "Large"
$ws.Range("D781").Value = "Center"
$ws.Range("E781").Value = "LTR"
$ws.Range("F781").Value = "<value>"

$ws.Range("B782").Value = "SingleUseId894"
$ws.Range("C782").Value = "Large"
$ws.Range("D782").Value = "Left"
$ws.Range("E782").Value = "LTR"
$ws.Range("F703").Copy() | Out-Null
$ws.Range("F782").PasteSpecial() | Out-Null

$excel.CutCopyMode = $false
